$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark bug in row 16 ("in create_word_xml, if the path arguments don't end
# with "/" then it fails...") as resolved, using same date format/style as
# other "Resolved" date cells (copy number format from D15).
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 43405

# Mark bug in row 23 ("need to make the xls file ignore populations if they
# are blank.") as resolved the same way.
$ws.Range("D15").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 43405

# Move the active selection to C16.
$ws.Range("C16").Select()
